$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 header values ----
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "setProperty"
$ws.Range("C1").Value = "open"
$ws.Range("D1").Value = "wait"
$ws.Range("E1").Value = "open"

# ---- Row 2 values ----
$ws.Range("B2").Value = '{"type":"json"}'
$ws.Range("C2").Value = "https://task.hugang.io/login"
$ws.Range("E2").Value = '${url}'

# ---- Row 3 values ----
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = '{"url":"https://task.hugang.io/"}'
$ws.Range("D3").Value = 2000

# ---- Apply the existing bordered/centered header style (same as A1/B1) to
#      every new cell that needs it, by copying B1's format (this re-uses
#      the existing style entry instead of minting new ones). ----
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Turn C2 into a hyperlink (adds the hyperlink font on top of the
#      border/center format already applied above). ----
$ws.Hyperlinks.Add($ws.Range("C2"), "https://task.hugang.io/login")

# ---- Column widths ----
$ws.Columns("B").ColumnWidth = 30.875
$ws.Columns("C").ColumnWidth = 27.5
$ws.Columns("D").ColumnWidth = 5.5
$ws.Columns("E").ColumnWidth = 6.375

# ---- Selection ----
$ws.Range("G4").Select()
